# The "M:" rule-type row (row 3 - "M:" / "Magic value program uses for some
# purpose internally") is no longer reported by the finance app, so it is
# removed from the legend table. Deleting the entire row shifts every row
# below it up by one, which also naturally drops the now-unused shared
# strings ("M:" and "Magic value program uses for some purpose internally")
# from the workbook and keeps every other row's styling/content intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

# Reflect where the user's cursor ended up after doing the edit.
$ws.Range("B10").Select()
